$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet2" to "Sheet 1"
$ws.Name = "Sheet 1"

# Rename "Persistence" -> "Permutation" everywhere relevant
$ws.Range("E17").Value = "Permutation imp."
$ws.Range("G17").Value = "Permutation"
$ws.Range("A13").Value = "Feature Selection: Using feature and permutation importance"

# Update the Table1 column header too (keeps in sync with E17, but set explicitly)
$wb.Sheets.Item("Sheet 1").ListObjects.Item("Table1").ListColumns.Item(5).Name = "Permutation imp."

# Update the selection / view to match the new state
$ws.Range("A13:E13").Select()
$ws.Application.ActiveWindow.ScrollRow = 13
